$d = $word.ActiveDocument

# The heading text "Example 4:  Variability of Athletes' Strengths - Residual
# Standard Deviation" needs to become "Example 14:  ...". This is the only
# visible textual change described by the diff (the other hunks are
# cosmetic run-splits/proofErr churn that leave the rendered text identical).

$d.Content.Find.Execute("Example 4:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Example 14:", 2)
